$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header-like values, extend with D1 and E1, copying the style from C1
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

# Row 2
$ws.Range("C2").Value = -4.993368022640259
$ws.Range("D2").Value = -4.657836315545413
$ws.Range("E2").Value = -4.29002747715822

# Row 3
$ws.Range("C3").Value = -1.165564360090414
$ws.Range("D3").Value = -1.161884795197474
$ws.Range("E3").Value = -1.137909370494665

# Row 4
$ws.Range("C4").Value = -0.06515441686030865
$ws.Range("D4").Value = 0.02779144705041408
$ws.Range("E4").Value = 0.09946321176602314

# Row 5
$ws.Range("C5").Value = -0.4032515873081615
$ws.Range("D5").Value = -0.1414060152494321
$ws.Range("E5").Value = 0.06543112287544997

# Row 6
$ws.Range("C6").Value = 0.01465567179956126
$ws.Range("D6").Value = -0.04421234966728196
$ws.Range("E6").Value = -0.09831220414754076

# Row 7
$ws.Range("C7").Value = 0.1084370207011733
$ws.Range("D7").Value = 0.06892470898996189
$ws.Range("E7").Value = 0.03132080520059642

# Row 8
$ws.Range("C8").Value = 0.1341971137761105
$ws.Range("D8").Value = 0.03973246936523338
$ws.Range("E8").Value = -0.04956464123919342
